$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.845.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.624.70'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.28'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.50%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.73'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.628.92'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.17'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.841.53'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.92'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0719'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.07'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.29'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.87'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0480'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.392.94'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.02'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +10.77%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.556'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.63'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.86%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.98'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.20%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.59'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.00%  '
